{"js": "// Insert three new paragraphs (a name/date header line, a \"Classe\" line,\n// and a trailing empty Title-styled paragraph) at the very beginning of\n// the document body, ahead of the existing \"Rapport du troubleshooting\n// DNS\" title paragraph.\n//\n// We build the content as a minimal Flat-OPC WordprocessingML package and\n// insert it with Body.insertOoxml(..., \"Start\") so we get exact control\n// over run/paragraph formatting (font size 10pt/11pt via sz/szCs, the\n// proofErr spell-check markers, and the literal tab runs) instead of the\n// formatting Word would otherwise inherit from the following Title\n// paragraph (centered, 22pt).\n//\n// Office.js merges the *last* paragraph of inserted OOXML into the\n// destination paragraph (same behavior as pasting), so a trailing empty\n// <w:p/> is appended to the package purely to absorb that merge and keep\n// the three authored paragraphs intact and separate.\nconst body = context.document.body;\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr>' +\n  '<w:pStyle w:val=\"Titre\"/>' +\n  '<w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr>' +\n  '</w:pPr>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t>Mayala-Luneko</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t xml:space=\"preserve\"> Loyde                                  </w:t></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:tab/></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:tab/></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:tab/></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:tab/></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:tab/></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:tab/></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:tab/></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:tab/><w:t>17/12/2024</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">Classe : </w:t></w:r>' +\n  '<w:r><w:t>2L1</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n  '<w:pPr>' +\n  '<w:pStyle w:val=\"Titre\"/>' +\n  '<w:rPr><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>' +\n  '</w:pPr>' +\n  '</w:p>' +\n  '<w:p/>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nbody.insertOoxml(flatOpcXml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Insert three new paragraphs (a name/date header line, a \"Classe\" line,\n# and a trailing empty Title-styled paragraph) at the very beginning of\n# the document, ahead of the existing \"Rapport du troubleshooting DNS\"\n# title paragraph.\n#\n# We build the content as a minimal Flat-OPC WordprocessingML package and\n# insert it at a zero-length Range at the start of the document via\n# Range.InsertXML so we get exact control over run/paragraph formatting\n# (font size 10pt/11pt via sz/szCs, the proofErr spell-check markers, and\n# the literal tab runs) instead of the formatting Word would otherwise\n# inherit from the following Title paragraph (centered, 22pt).\n\n$d = $word.ActiveDocument\n\n$flatOpcXml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p>\n  <w:pPr>\n    <w:pStyle w:val=\"Titre\"/>\n    <w:rPr>\n      <w:sz w:val=\"20\"/>\n      <w:szCs w:val=\"20\"/>\n    </w:rPr>\n  </w:pPr>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"20\"/>\n      <w:szCs w:val=\"20\"/>\n    </w:rPr>\n    <w:t>Mayala-Luneko</w:t>\n  </w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"20\"/>\n      <w:szCs w:val=\"20\"/>\n    </w:rPr>\n    <w:t xml:space=\"preserve\"> Loyde                                  </w:t>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"20\"/>\n      <w:szCs w:val=\"20\"/>\n    </w:rPr>\n    <w:tab/>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"20\"/>\n      <w:szCs w:val=\"20\"/>\n    </w:rPr>\n    <w:tab/>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"20\"/>\n      <w:szCs w:val=\"20\"/>\n    </w:rPr>\n    <w:tab/>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"20\"/>\n      <w:szCs w:val=\"20\"/>\n    </w:rPr>\n    <w:tab/>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"20\"/>\n      <w:szCs w:val=\"20\"/>\n    </w:rPr>\n    <w:tab/>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"20\"/>\n      <w:szCs w:val=\"20\"/>\n    </w:rPr>\n    <w:tab/>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"20\"/>\n      <w:szCs w:val=\"20\"/>\n    </w:rPr>\n    <w:tab/>\n  </w:r>\n  <w:r>\n    <w:rPr>\n      <w:sz w:val=\"20\"/>\n      <w:szCs w:val=\"20\"/>\n    </w:rPr>\n    <w:tab/>\n    <w:t>17/12/2024</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:r>\n    <w:t xml:space=\"preserve\">Classe : </w:t>\n  </w:r>\n  <w:r>\n    <w:t>2L1</w:t>\n  </w:r>\n</w:p>\n<w:p>\n  <w:pPr>\n    <w:pStyle w:val=\"Titre\"/>\n    <w:rPr>\n      <w:sz w:val=\"22\"/>\n      <w:szCs w:val=\"22\"/>\n    </w:rPr>\n  </w:pPr>\n</w:p>\n</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n\n$r = $d.Range(0, 0)\n$r.InsertXML($flatOpcXml)\n"}
